$d = $word.ActiveDocument

# Locate the "Here is a photo of Cecil:" caption paragraph and the following
# paragraph that holds the inline picture, then remove both (text, runs,
# drawing and paragraph marks) in one go.
$captionPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Here is a photo of Cecil:*") {
        $captionPara = $p
        break
    }
}

if ($captionPara -ne $null) {
    $captionIndex = $captionPara.Index
    $photoPara = $d.Paragraphs.Item($captionIndex + 1)
    $deleteRange = $d.Range($captionPara.Range.Start, $photoPara.Range.End)
    $deleteRange.Delete()
}

# Move the hidden "_GoBack" bookmark from the end of the "Cecil is ten years
# old ..." paragraph to the very start of that same paragraph.
$birthdayPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Cecil is ten years old*") {
        $birthdayPara = $p
        break
    }
}

if ($birthdayPara -ne $null) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
    $startPoint = $d.Range($birthdayPara.Range.Start, $birthdayPara.Range.Start)
    $d.Bookmarks.Add("_GoBack", $startPoint)
}
